$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Add the new entry as row 78
$ws.Range("A78").Value = 45405
$ws.Range("B78").Value = 6
$ws.Range("C78").Value = "QSE Bugfixes bei QueryBased, Vergleiche zwischen QB und FB"

# Match the style used for the date column above (format code 16, e.g. "d-mmm")
$ws.Range("A78").NumberFormat = $ws.Range("A77").NumberFormat

# Update the active selection to the newly-added cell
$ws.Range("C78").Select()
